$d = $word.ActiveDocument

# 1. Collapse the title run-split (caused by spell-check proofErr wrapping) into one run.
$d.Content.Find.Execute("Aanpassingen applicatie n.a.v testbevindigen", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Aanpassingen applicatie n.a.v testbevindigen", 2)

# 2. Update the date day-of-month from 21 to 24.
$d.Content.Find.Execute("Datum: 21-05-2018", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Datum: 24-05-2018", 2)

$d.Content.SpellingChecked = $true
